$wb = $excel.ActiveWorkbook

# Fix the typo'd "maxDetxerity" header -> "maxDexterity" on every class sheet.
# Replacing every reference lets the writer drop the now-unused shared
# string and append the corrected one at the end of the table.
foreach ($name in @("warrior", "mage", "rogue", "hunter", "merchant")) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("H1").Value = "maxDexterity"
}

# Restore per-sheet selections (order matters: the last sheet touched
# becomes the active/tabSelected sheet, and "merchant" must stay active).
$wsWarrior = $wb.Worksheets.Item("warrior")
$wsWarrior.Range("G10").Select()

$wsMage = $wb.Worksheets.Item("mage")
$wsMage.Range("D37").Select()

$wsRogue = $wb.Worksheets.Item("rogue")
$wsRogue.Range("H7").Select()

$wsHunter = $wb.Worksheets.Item("hunter")
$wsHunter.Range("H9").Select()

$wsMerchant = $wb.Worksheets.Item("merchant")
$wsMerchant.Range("I9").Select()
$excel.ActiveWindow.Zoom = 85
